$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''29.955.15'
$ws.Range("E2").Value = '  +0.45%  '
$ws.Range("D3").Value = '''1.893.46'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''0.7770'
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").Value = '''244.01'
$ws.Range("E6").Value = '  -0.03%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '''0.3128'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").Value = '''25.86'
$ws.Range("E9").Value = '  +1.72%  '
$ws.Range("D10").Value = '''0.07248'
$ws.Range("E10").Value = '  +0.70%  '
$ws.Range("D11").Value = '''0.08700'
$ws.Range("E11").Value = '  +7.69%  '
$ws.Range("D12").Value = '''2.102.18'
$ws.Range("E12").Value = '  +10.20%  '
$ws.Range("D13").Value = '''0.7748'
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("D14").Value = '''5.417'
$ws.Range("E14").Value = '  -1.35%  '
$ws.Range("D15").Value = '''94.53'
$ws.Range("E15").Value = '  +2.27%  '
$ws.Range("D16").Value = '''6.184'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").Value = '''30.189.51'
$ws.Range("E17").Value = '  +1.12%  '
$ws.Range("D18").Value = '''13.90'
$ws.Range("E18").Value = '  -0.50%  '
$ws.Range("B19").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C19").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D19").Value = '''2.333.07'
$ws.Range("E19").Value = '  +5.71%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").Value = '''246.14'
$ws.Range("E20").Value = '  +0.90%  '
$ws.Range("D21").Value = '''0.000007865'
$ws.Range("E21").Value = '  +1.23%  '
$ws.Range("B22").Value = 'Dai'
$ws.Range("C22").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D22").Value = '''1.001'
$ws.Range("E22").Value = '  +0.16%  '
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").Value = '''8.136'
$ws.Range("E23").Value = '  -0.25%  '
$ws.Range("E24").Value = '  +0.04%  '
$ws.Range("D25").Value = '''0.1661'
$ws.Range("E25").Value = '  +6.84%  '
$ws.Range("D26").Value = '''9.489'
$ws.Range("E26").Value = '  +0.91%  '
$ws.Range("D27").Value = '''163.34'
$ws.Range("E27").Value = '  +0.47%  '
$ws.Range("E28").Value = '  +0.51%  '
$ws.Range("D29").Value = '''2.054'
$ws.Range("E29").Value = '  +0.21%  '
$ws.Range("D30").Value = '''1.434'
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("D32").Value = '''4.520'
$ws.Range("E32").Value = '  +1.00%  '
$ws.Range("D33").Value = '''4.132'
$ws.Range("E33").Value = '  +0.59%  '
$ws.Range("D34").Value = '''0.05475'
$ws.Range("E34").Value = '  -1.01%  '
$ws.Range("D35").Value = '''1.245'
$ws.Range("E35").Value = '  -1.43%  '
$ws.Range("D36").Value = '''0.7559'
$ws.Range("E37").Value = '  +0.14%  '
$ws.Range("D38").Value = '''2.691'
$ws.Range("E38").Value = '  +2.39%  '
$ws.Range("E39").Value = '  +2.52%  '
$ws.Range("D40").Value = '''2.788'
$ws.Range("E40").Value = '  +0.36%  '
$ws.Range("E41").Value = '  +1.82%  '
$ws.Range("B42").Value = 'Maker'
$ws.Range("C42").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D42").Value = '''1.109.03'
$ws.Range("E42").Value = '  -2.65%  '
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '''73.77'
$ws.Range("E43").Value = '  +0.01%  '
$ws.Range("D44").Value = '''6.103'
$ws.Range("E44").Value = '  +3.42%  '
$ws.Range("D45").Value = '''0.8538'
$ws.Range("E45").Value = '  +0.31%  '
$ws.Range("D46").Value = '''2.194.26'
$ws.Range("E46").Value = '  +6.32%  '
$ws.Range("E47").Value = '  +0.09%  '
$ws.Range("D48").Value = '''103.66'
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("E49").Value = '  -0.72%  '
$ws.Range("D50").Value = '''7.611'
$ws.Range("E50").Value = '  +1.89%  '
$ws.Range("D51").Value = '''9.890'
$ws.Range("E51").Value = '  -0.35%  '
